# Update the actuals for the "update db via react app" task (row 19):
#  - ActlHours (P19):        6  -> 9
#  - ActDeliveryDate (Q19):  2020-02-11 Pending -> 2020-02-14 Pending
#  - Weekday (R19):          Tuesday -> Friday   (2020-02-14 is a Friday)
#  - Comment (T19) stays the same text, just re-select so the active cell
#    matches the saved view (R19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P19").Value = 9
$ws.Range("Q19").Value = "2020-02-14 Pending"
$ws.Range("R19").Value = "Friday"
$ws.Range("T19").Value = "Greate trouble making it work to update db via react app"

$ws.Range("R19").Select()
